$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells: force Text format so numeric-looking strings
# (e.g. "548.82", "0.500", multi-dot "63.469.13") are kept as literal text
# instead of being auto-coerced into floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.469.13"
$ws.Range("E2").Value = "  +5.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.063.27"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.82"
$ws.Range("E5").Value = "  +5.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.83"
$ws.Range("E6").Value = "  +7.87%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.058.36"
$ws.Range("E8").Value = "  +3.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").Value = "  +3.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.23"
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  +2.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  +4.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("E13").Value = "  +5.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.65"
$ws.Range("E14").Value = "  +5.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.569.95"
$ws.Range("E15").Value = "  +3.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.577.82"
$ws.Range("E16").Value = "  +5.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.069.95"
$ws.Range("E17").Value = "  +3.88%  "
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.73"
$ws.Range("E19").Value = "  +5.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.49"
$ws.Range("E20").Value = "  +6.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.56"
$ws.Range("E21").Value = "  +5.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("E22").Value = "  +2.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.20"
$ws.Range("E23").Value = "  +7.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.48"
$ws.Range("E24").Value = "  +5.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.52"
$ws.Range("E25").Value = "  +7.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("E27").Value = "  +5.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.93"
$ws.Range("E28").Value = "  +4.60%  "
$ws.Range("E29").Value = "  +9.55%  "
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.97"
$ws.Range("E31").Value = "  +4.55%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.42"
$ws.Range("E33").Value = "  +9.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.67"
$ws.Range("E34").Value = "  +7.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.79"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.98"
$ws.Range("E36").Value = "  +5.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "467.07"
$ws.Range("E37").Value = "  +4.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0815"
$ws.Range("E38").Value = "  +6.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.143.92"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0396"
$ws.Range("E40").Value = "  +6.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.24"
$ws.Range("E42").Value = "  +4.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("E43").Value = "  +8.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.97"
$ws.Range("E44").Value = "  +12.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.252"
$ws.Range("E45").Value = "  +5.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.05"
$ws.Range("E47").Value = "  +7.59%  "
$ws.Range("E48").Value = "  +2.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0511"
$ws.Range("E49").Value = "  +3.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "116.19"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("E51").Value = "  +7.60%  "
